# Fix a bug in PDS: add three missing "Ignore" filter entries
# (electricity meter / water meter / gas meter) to the Filter sheet,
# inserted right before the existing E-BFAS630 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filter")

# Insert three new blank rows above the current row 18 ("E-BFAS630" ...)
$ws.Rows("18:20").Insert()

# Row 18: E-BEMS100 / Ignore / 电能表 (electricity meter)
$ws.Range("A18").Value = "E-BEMS100"
$ws.Range("C18").Value = "Ignore"
$ws.Range("D18").Value = "电能表"

# Row 19: E-BEMS110 / Ignore / 水表 (water meter)
$ws.Range("A19").Value = "E-BEMS110"
$ws.Range("C19").Value = "Ignore"
$ws.Range("D19").Value = "水表"

# Row 20: E-BEMS120 / Ignore / 燃气表 (gas meter)
$ws.Range("A20").Value = "E-BEMS120"
$ws.Range("C20").Value = "Ignore"
$ws.Range("D20").Value = "燃气表"

# Leave the active selection on the first newly-added row, as in the source file
$ws.Range("A18").Select()
